$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new rows of glucose readings (rows 30 and 31), continuing the
# existing daily log. Formats are inherited from the row directly above
# (row 29), matching Excel's default behavior when extending a table.

$ws.Range("A30").Value = 44223
$ws.Range("B30").Value = 4.4000000000000004
$ws.Range("C30").Value = 6.5
$ws.Range("D30").Value = 4.7
$ws.Range("E30").Value = 5.6
$ws.Range("F30").Value = 4.5999999999999996
$ws.Range("G30").Value = 4.3

$ws.Range("A31").Value = 44224
$ws.Range("B31").Value = 4.4000000000000004
$ws.Range("C31").Value = 5.4
$ws.Range("D31").Value = 5.2
$ws.Range("E31").Value = 5.3
$ws.Range("F31").Value = 4.9000000000000004
$ws.Range("G31").Value = 4.9000000000000004

# Keep the same formatting as the row above (row 29) by copying its
# formats down, the same way Excel extends an existing table when you
# type into the row right below it.
$ws.Range("A29:G29").Copy()
$ws.Range("A30:G31").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the active selection to match what was left selected after typing
# the new data (A30:A31).
$ws.Range("A30:A31").Select()
